# Daily attendance processing - 2025-10-28 23:42:55
# Normalize the "Recorded By" (column G) author-list ordering so the
# primary human/system recorder is listed ahead of the generic "System"
# token (and the lower-case "system" dedupe token sorts after "System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every row in this report (header row 1 excluded) that has a value in
# column G is a candidate; only the two known mis-ordered patterns are
# rewritten, everything else is left untouched.
$targetRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25, `
  29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52, `
  56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79, `
  83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103, `
  109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129, `
  135,136,137,138,139,140,141,142,143,144,145,146,147,148,149,150,151,152,153,154,155)

foreach ($r in $targetRows) {
  $cell = $ws.Range("G$r")
  $old = $cell.Value()
  if ($old -eq "System, dnasr281@gmail.com") {
    $cell.Value = "dnasr281@gmail.com, System"
  } elseif ($old -eq "backup@backdoor.com, system, System") {
    $cell.Value = "backup@backdoor.com, System, system"
  }
}

Write-Host "Recorded By ordering normalized"
